# Remove the replicate microstate "SM24_micro008" from the microstate list.
#
# The sheet lists one microstate per row (columns B/C = ID, SMILES), each
# paired with a 2D-depiction picture anchored to that row. Row 9 holds the
# replicate entry SM24_micro008. Removing it shifts every following
# molecule's ID/SMILES up by one row (values only - the per-row banding
# style stays put), drops the now-surplus last row (29) and the
# now-surplus last picture in the drawing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 29
$idCol = 2      # column B: microstate ID
$smilesCol = 3  # column C: canonical isomeric SMILES

# Shift the ID/SMILES values up by one row, starting at the replicate's
# row (9), overwriting it with what followed. This only moves the cell
# VALUES, leaving each row's existing formatting (banding) untouched.
for ($r = 9; $r -lt $lastRow; $r++) {
    $ws.Cells.Item($r, $idCol).Value = $ws.Cells.Item($r + 1, $idCol).Value()
    $ws.Cells.Item($r, $smilesCol).Value = $ws.Cells.Item($r + 1, $smilesCol).Value()
}

# The data now ends one row earlier - drop the trailing (now duplicated) row.
$ws.Rows($lastRow).Delete()

# Drop the now-surplus last 2D-depiction picture to match the new row count.
$ws.Shapes.Item($ws.Shapes.Count).Delete()
